$d = $word.ActiveDocument

# --- Part 1: split the "===> Resulting in 2^4 = 16 ... different spaces to
#     clean" paragraph into two paragraphs, replacing the leading "===> "
#     marker with a Wingdings arrow symbol, and moving the _GoBack bookmark
#     into a new second paragraph that introduces the room-selection
#     conclusion. ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "===>*different spaces to clean*") {
        $target = $p.Range
    }
}

$xmlPart1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:sym w:font="Wingdings" w:char="F0E8"/></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t>Resulting in 2^4 = 16</w:t></w:r>
<w:r><w:t xml:space="preserve"> different spaces to clean</w:t></w:r>
</w:p>
<w:p>
<w:r><w:sym w:font="Wingdings" w:char="F0E8"/></w:r>
<w:r><w:t xml:space="preserve"> Conclusion on what rooms to use: </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlPart1)

# --- Part 2: append a new "06-05-2022, Friday week 2" heading paragraph
#     (same Kop1 style as the other date headings) plus a trailing empty
#     paragraph at the very end of the document. ---
$endPos = $d.Content.End
$tail = $d.Range($endPos, $endPos)

$xmlPart2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:pStyle w:val="Kop1"/></w:pPr>
<w:r><w:t>06</w:t></w:r>
<w:r><w:t xml:space="preserve">-05-2022, </w:t></w:r>
<w:r><w:t>Friday</w:t></w:r>
<w:r><w:t xml:space="preserve"> week 2</w:t></w:r>
</w:p>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$tail.InsertXML($xmlPart2)
